$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: update AF6 (related_works)
$ws.Range("AF6").Value = 'c("https://openalex.org/W2953102189", "https://openalex.org/W2618406201", "https://openalex.org/W2565995904", "https://openalex.org/W2293197596", "https://openalex.org/W2998674080", "https://openalex.org/W2281217800", "https://openalex.org/W2766491241", "https://openalex.org/W155792083", "https://openalex.org/W1519874714", "https://openalex.org/W1968964438")'

# Row 8: update F8 (so), G8 (so_id), H8 (host_organization), I8 (issn_l)
$ws.Range("F8").Value = "Clinical Lymphoma, Myeloma & Leukemia"
$ws.Range("G8").Value = "https://openalex.org/S186546955"
$ws.Range("H8").Value = "Elsevier BV"
$ws.Range("I8").Value = "2152-2669"

# Row 9: update F9 (so), G9 (so_id), H9 (host_organization), I9 (issn_l)
$ws.Range("F9").Value = "Avicenna Journal of Medicine"
$ws.Range("G9").Value = "https://openalex.org/S2764417561"
$ws.Range("H9").Value = "Thieme Medical Publishers (Germany)"
$ws.Range("I9").Value = "2231-0770"

# Row 9: update AF9 (related_works)
$ws.Range("AF9").Value = 'c("https://openalex.org/W2096401073", "https://openalex.org/W3146331961", "https://openalex.org/W2033331561", "https://openalex.org/W3114922805", "https://openalex.org/W2419073828", "https://openalex.org/W2792196017", "https://openalex.org/W1970421407", "https://openalex.org/W170943972", "https://openalex.org/W2047981850", "https://openalex.org/W2051555008")'
